# Add "I0" and "IF" columns (I and J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy H1's formatting (bold, border, centered) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: row number -> (I value, J value)
$data = @{
    2  = @(2, 6)
    3  = @(1, 7)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(1, 7)
    7  = @(1, 7)
    8  = @(1, 7)
    9  = @(1, 5)
    10 = @(1, 8)
    11 = @(1, 6)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(1, 4)
    15 = @(1, 5)
    16 = @(1, 6)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 5)
    20 = @(1, 6)
    21 = @(1, 5)
    22 = @(5, 9)
    23 = @(1, 4)
    24 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
